$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 1 values (columns B..K) ---
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 21
$ws.Range("D1").Value = 32
$ws.Range("E1").Value = 26
$ws.Range("F1").Value = 20
$ws.Range("G1").Value = 18
$ws.Range("H1").Value = 10
$ws.Range("I1").Value = 0.069999999999999993
$ws.Range("J1").Value = 0.017000000000000001
$ws.Range("K1").Value = 0.098999999999999991

# --- Update column widths for columns F, H, I, J (6, 8, 9, 10) ---
# The OOXML "width" attribute is stored in character-width units derived from
# pixel widths; ColumnWidth here snaps to the nearest achievable pixel value,
# so we pick the ColumnWidth value whose resulting pixel width is closest to
# the target stored width.
$ws.Columns.Item(6).ColumnWidth = 2.3333333333333335   # -> stored width ~3.140625 (target col F)
$ws.Columns.Item(8).ColumnWidth = 2.3333333333333335   # -> stored width ~3.140625 (target col H)
$ws.Columns.Item(9).ColumnWidth = 3.8333333333333335   # -> stored width ~4.7109375 (target col I)
$ws.Columns.Item(10).ColumnWidth = 4.833333333333333   # -> stored width ~5.7109375 (target col J)
